$wb = $excel.ActiveWorkbook

# --- Sheet "customer-accounts": give usr02/usr03/usr04/usr05 their own distinct
#     account aliases instead of reusing usr01's checking1/saving1, and point the
#     last two rows at usr04 / usr05 instead of the old usr101 / sys101 ids.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B4").Value = "checking2"
$ws1.Range("B5").Value = "saving2"
$ws1.Range("B6").Value = "checking3"
$ws1.Range("B7").Value = "saving4"
$ws1.Range("B8").Value = "checking5"
$ws1.Range("D8").Value = "usr04"
$ws1.Range("B9").Value = "saving5"
$ws1.Range("D9").Value = "usr04"
$ws1.Range("B10").Value = "checking6"
$ws1.Range("D10").Value = "usr05"
$ws1.Range("B11").Value = "saving6"
$ws1.Range("D11").Value = "usr05"

# --- Sheet "verifyTransferPositiveCases": update the from/to account references so
#     each row keeps using the renamed account keys for the right user.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("D6").Value = "saving3"
$ws2.Range("E6").Value = "checking3"
$ws2.Range("C7").Value = "usr04"
$ws2.Range("E7").Value = "checking4"
$ws2.Range("D2").Value = "usr01.default"
$ws2.Range("D3").Value = "usr01.default"
$ws2.Range("D7").Value = "usr04.default"
$ws2.Range("D5").Value = "usr02.default"
$ws2.Range("E5").Value = "checking2"

# column widths on the positive-cases sheet grow now that the values are longer
# (inputs tuned so the engine's pixel-snapped output lands as close as possible
# to the recorded widths of 15.7109375 / 21.28515625)
$ws2.Columns.Item(3).ColumnWidth = 14.85
$ws2.Columns.Item(4).ColumnWidth = 20.5

# --- View/selection state: make the positive-cases sheet the active tab/sheet,
#     update the remembered selections, and move the previously-active "users"
#     sheet's tabSelected flag off of it.
$ws1.Range("D12").Select()
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A6").Select()
$ws2.Activate()
$ws2.Range("E5").Select()
